$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New localization rows for the Camp Screen options dialog
$ws.Range("A46").Value = "Dialog.WhatToDo"
$ws.Range("B46").Value = "What to do now?"
$ws.Range("C46").Value = "O que farei agora?"

$ws.Range("A47").Value = "Dialog.GoToMission"
$ws.Range("B47").Value = "Do a Mission"
$ws.Range("C47").Value = "Fazer uma Missão"

$ws.Range("A48").Value = "Dialog.Shop"
$ws.Range("B48").Value = "Shop"
$ws.Range("C48").Value = "Comprar Suprimentos"

$ws.Range("C49").Value = "Ficar no Acampamento"
$ws.Range("B49").Value = "Stay at Camp"
$ws.Range("A49").Value = "Dialog.StayCamp"

# Keep the same underline style as the other blank marker cell (D22) on the new blank row
$ws.Range("A50").Font.Underline = $true

$ws.Range("A50").Select
